# Updated symbol list on Mon Jan 30 15:13:32 UTC 2023 with GitHub Actions
#
# Refreshes the crypto price table (columns D=Price, E=Volume(1h), G=Hora)
# for every data row (2-51) on the active sheet with the latest scrape.
# Values are written as plain text (matching the sheet's existing
# inline-string cells) rather than numbers/percentages, so each cell is
# forced to the "Text" number format before assignment and restored to
# the default "Normal" style afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row number -> column letter -> new text value
$updates = @{
    2 = @{ "D" = "309.52"; "E" = "-2.84%"; "G" = "15" }
    3 = @{ "D" = "38.05"; "E" = "-3.61%"; "G" = "15" }
    4 = @{ "D" = "5.076"; "E" = "-0.99%"; "G" = "15" }
    5 = @{ "D" = "0.07892"; "E" = "-3.72%"; "G" = "15" }
    6 = @{ "E" = "-3.11%"; "G" = "15" }
    7 = @{ "D" = "4.361"; "E" = "1.88%"; "G" = "15" }
    8 = @{ "D" = "8.312"; "E" = "0.50%"; "G" = "15" }
    9 = @{ "D" = "3.035"; "E" = "-7.45%"; "G" = "15" }
    10 = @{ "D" = "0.9308"; "E" = "-0.52%"; "G" = "15" }
    11 = @{ "D" = "0.1300"; "E" = "-7.58%"; "G" = "15" }
    12 = @{ "D" = "0.1946"; "E" = "-2.37%"; "G" = "15" }
    13 = @{ "D" = "0.08776"; "E" = "-3.29%"; "G" = "15" }
    14 = @{ "D" = "0.03414"; "E" = "-4.43%"; "G" = "15" }
    15 = @{ "D" = "0.09743"; "E" = "-0.69%"; "G" = "15" }
    16 = @{ "D" = "0.001389"; "E" = "-0.79%"; "G" = "15" }
    17 = @{ "D" = "0.005942"; "E" = "-1.75%"; "G" = "15" }
    18 = @{ "E" = "1,777.17%"; "G" = "15" }
    19 = @{ "D" = "3.593"; "E" = "-1.94%"; "G" = "15" }
    20 = @{ "D" = "0.3436"; "E" = "-0.74%"; "G" = "15" }
    21 = @{ "D" = "0.1295"; "E" = "-0.65%"; "G" = "15" }
    22 = @{ "D" = "4.992"; "E" = "1.72%"; "G" = "15" }
    23 = @{ "D" = "0.2484"; "E" = "1.36%"; "G" = "15" }
    24 = @{ "D" = "0.04309"; "E" = "-0.41%"; "G" = "15" }
    25 = @{ "D" = "0.001216"; "E" = "-0.75%"; "G" = "15" }
    26 = @{ "E" = "-3.45%"; "G" = "15" }
    27 = @{ "E" = "176.48%"; "G" = "15" }
    28 = @{ "G" = "15" }
    29 = @{ "G" = "15" }
    30 = @{ "G" = "15" }
    31 = @{ "G" = "15" }
    32 = @{ "G" = "15" }
    33 = @{ "G" = "15" }
    34 = @{ "G" = "15" }
    35 = @{ "G" = "15" }
    36 = @{ "G" = "15" }
    37 = @{ "G" = "15" }
    38 = @{ "G" = "15" }
    39 = @{ "D" = "0.02303"; "E" = "3.45%"; "G" = "15" }
    40 = @{ "D" = "0.05052"; "E" = "-3.90%"; "G" = "15" }
    41 = @{ "D" = "0.007498"; "E" = "-0.20%"; "G" = "15" }
    42 = @{ "D" = "0.009903"; "E" = "0.00%"; "G" = "15" }
    43 = @{ "D" = "0.1360"; "E" = "-1.47%"; "G" = "15" }
    44 = @{ "D" = "0.002017"; "E" = "-6.15%"; "G" = "15" }
    45 = @{ "D" = "0.008765"; "E" = "-11.16%"; "G" = "15" }
    46 = @{ "D" = "0.00006553"; "E" = "-0.94%"; "G" = "15" }
    47 = @{ "D" = "0.00000000749"; "E" = "-0.10%"; "G" = "15" }
    48 = @{ "D" = "0.002993"; "E" = "8.09%"; "G" = "15" }
    49 = @{ "G" = "15" }
    50 = @{ "D" = "0.00002097"; "E" = "-0.10%"; "G" = "15" }
    51 = @{ "D" = "0.0001997"; "E" = "-0.10%"; "G" = "15" }
}

foreach ($row in $updates.Keys) {
    $rowData = $updates[$row]
    foreach ($col in $rowData.Keys) {
        $addr = "$col$row"
        $value = $rowData[$col]
        $cell = $ws.Range($addr)
        $cell.NumberFormat = "@"
        $cell.Value = $value
        $cell.Style = "Normal"
    }
}
